$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player rows (20-31) continue the existing 3-row-per-player block,
# alternating the two fill styles already used in the table (green / yellow).
# Copy the formatting from an existing same-style block so the saved file
# reuses the same cellXfs/fill entries rather than creating new ones.
$xlPasteFormats = -4122

$ws.Range("A2:F4").Copy()
$ws.Range("A20:F22").PasteSpecial($xlPasteFormats)
$ws.Range("A5:F7").Copy()
$ws.Range("A23:F25").PasteSpecial($xlPasteFormats)
$ws.Range("A2:F4").Copy()
$ws.Range("A26:F28").PasteSpecial($xlPasteFormats)
$ws.Range("A5:F7").Copy()
$ws.Range("A29:F31").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Row 20: Harrison Smith / Group1
$ws.Cells.Item(20,1).Value = "Harrison Smith"
$ws.Cells.Item(20,2).Value = "Group1"
$ws.Cells.Item(20,3).Value = 9.333333333333334
$ws.Cells.Item(20,4).Value = 96
$ws.Cells.Item(20,5).Value = 67.33333333333333
$ws.Cells.Item(20,6).Value = 28.66666666666667

# Row 21: Harrison Smith / Group2
$ws.Cells.Item(21,1).Value = "Harrison Smith"
$ws.Cells.Item(21,2).Value = "Group2"
$ws.Cells.Item(21,3).Value = 7.666666666666667
$ws.Cells.Item(21,4).Value = 88.33333333333333
$ws.Cells.Item(21,5).Value = 60.66666666666666
$ws.Cells.Item(21,6).Value = 27.66666666666667

# Row 22: Harrison Smith / Difference
$ws.Cells.Item(22,1).Value = "Harrison Smith"
$ws.Cells.Item(22,2).Value = "Difference"
$ws.Cells.Item(22,3).Value = -1.666666666666667
$ws.Cells.Item(22,4).Value = -7.666666666666671
$ws.Cells.Item(22,5).Value = -6.666666666666664
$ws.Cells.Item(22,6).Value = -1

# Row 23: Jalen Mills / Group1
$ws.Cells.Item(23,1).Value = "Jalen Mills"
$ws.Cells.Item(23,2).Value = "Group1"
$ws.Cells.Item(23,3).Value = 5.666666666666667
$ws.Cells.Item(23,4).Value = 54
$ws.Cells.Item(23,5).Value = 41
$ws.Cells.Item(23,6).Value = 13

# Row 24: Jalen Mills / Group2
$ws.Cells.Item(24,1).Value = "Jalen Mills"
$ws.Cells.Item(24,2).Value = "Group2"
$ws.Cells.Item(24,3).Value = 4
$ws.Cells.Item(24,4).Value = 40
$ws.Cells.Item(24,5).Value = 23
$ws.Cells.Item(24,6).Value = 17

# Row 25: Jalen Mills / Difference
$ws.Cells.Item(25,1).Value = "Jalen Mills"
$ws.Cells.Item(25,2).Value = "Difference"
$ws.Cells.Item(25,3).Value = -1.666666666666667
$ws.Cells.Item(25,4).Value = -14
$ws.Cells.Item(25,5).Value = -18
$ws.Cells.Item(25,6).Value = 4

# Row 26: Julian Love / Group1
$ws.Cells.Item(26,1).Value = "Julian Love"
$ws.Cells.Item(26,2).Value = "Group1"
$ws.Cells.Item(26,3).Value = 4.333333333333333
$ws.Cells.Item(26,4).Value = 55.66666666666666
$ws.Cells.Item(26,5).Value = 37.66666666666666
$ws.Cells.Item(26,6).Value = 18

# Row 27: Julian Love / Group2
$ws.Cells.Item(27,1).Value = "Julian Love"
$ws.Cells.Item(27,2).Value = "Group2"
$ws.Cells.Item(27,3).Value = 9
$ws.Cells.Item(27,4).Value = 118.6666666666667
$ws.Cells.Item(27,5).Value = 81
$ws.Cells.Item(27,6).Value = 37.66666666666666

# Row 28: Julian Love / Difference
$ws.Cells.Item(28,1).Value = "Julian Love"
$ws.Cells.Item(28,2).Value = "Difference"
$ws.Cells.Item(28,3).Value = 4.666666666666667
$ws.Cells.Item(28,4).Value = 63.00000000000001
$ws.Cells.Item(28,5).Value = 43.33333333333334
$ws.Cells.Item(28,6).Value = 19.66666666666666

# Row 29: Marcus Williams / Group1
$ws.Cells.Item(29,1).Value = "Marcus Williams"
$ws.Cells.Item(29,2).Value = "Group1"
$ws.Cells.Item(29,3).Value = 9.333333333333334
$ws.Cells.Item(29,4).Value = 62.66666666666666
$ws.Cells.Item(29,5).Value = 43.66666666666666
$ws.Cells.Item(29,6).Value = 19

# Row 30: Marcus Williams / Group2
$ws.Cells.Item(30,1).Value = "Marcus Williams"
$ws.Cells.Item(30,2).Value = "Group2"
$ws.Cells.Item(30,3).Value = 6
$ws.Cells.Item(30,4).Value = 49.66666666666666
$ws.Cells.Item(30,5).Value = 35
$ws.Cells.Item(30,6).Value = 14.66666666666667

# Row 31: Marcus Williams / Difference
$ws.Cells.Item(31,1).Value = "Marcus Williams"
$ws.Cells.Item(31,2).Value = "Difference"
$ws.Cells.Item(31,3).Value = -3.333333333333334
$ws.Cells.Item(31,4).Value = -13
$ws.Cells.Item(31,5).Value = -8.666666666666664
$ws.Cells.Item(31,6).Value = -4.333333333333334

